$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row appended to the rate history table.
$row = 58

# Column A holds a literal date-like string (e.g. "2025-09-04") rather than a
# real Excel date, matching every other row in the sheet. Force the cell to
# Text first so assigning the value doesn't get auto-converted to a date
# serial number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-10-04"
$ws.Cells.Item($row, 2).Value = "21:21:10"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,794.1737"
